# Applies the "all source and test" commit: reorders/updates the existing
# test-case rows, adds three new Servlet test cases (CART_SERV_08, 07, 09),
# and adjusts the Sample-Data / Steps column widths accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (B: Scenario, C: Sample Data, D: Steps, E: Expected) ---
$ws.Columns.Item(2).ColumnWidth = 24.75
$ws.Columns.Item(3).ColumnWidth = 11.916666666666666
$ws.Columns.Item(4).ColumnWidth = 29.583333333333332
$ws.Columns.Item(5).ColumnWidth = 22.083333333333332

# --- Full target table (rows 2-10): ID, Scenario, Sample Data, Steps, Expected, Actual, Status ---
$data = @(
  @("CART_SERV_08", "ID Sản phẩm lỗi format", "ID='abc'", "ID='abc' -> Exception", "Vào catch -> Báo lỗi", "OK", "PASS"),
  @("CART_SERV_02", "Servlet: Thêm thành công", "PID: 1, Qty: 2", "1. User OK`n2. Mock DAO tìm thấy SP`n3. Service add OK", "Redirect: Referer", "OK", "PASS"),
  @("CART_SERV_05", "Servlet: Số lượng âm", "Qty='-5'", "Input quantity='-5'", "Vẫn gọi service", "OK", "PASS"),
  @("CART_SERV_04", "Servlet: Số lượng là chữ", "Qty='abc'", "Input quantity='abc'", "Mặc định thêm 1", "OK", "PASS"),
  @("CART_SERV_01", "Servlet: Chưa đăng nhập", "User: null", "1. User session = null`n2. Call doPost", "Redirect: Login.jsp", "OK", "PASS"),
  @("CART_SERV_07", "Không có Referer Header", "Ref=null", "Header Referer = null", "Redirect default path", "OK", "PASS"),
  @("CART_SERV_03", "Servlet: SP không tồn tại", "PID: 999", "1. User OK`n2. Mock DAO trả về null", "Báo lỗi addCartError", "OK", "PASS"),
  @("CART_SERV_09", "Giỏ hàng trong Session Null", "Cart=null", "Session.getAttribute('cart') = null", "Service được gọi với null", "OK", "PASS"),
  @("CART_SERV_06", "Servlet: Lỗi Database", "DB Error", "DAO ném Exception", "Redirect & Báo lỗi", "OK", "PASS")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}

# --- Re-apply the existing PASS (green bold) status style to every status cell ---
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G2:G10").PasteSpecial(-4122) | Out-Null
